# Recompute NATMI LR-pair (Vegfb-Nrp1) stats with the updated TPM inputs.
# Only the TPM-derived value columns (G,H,I,J,M,N,O,P,Q,R,S,T) change;
# the identifying columns (A-F,K,L) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 1.844441333333333  # G2: Ligand average expression value
$ws.Cells.Item(2, 8).Value2 = 5.533324  # H2: Ligand total expression value
$ws.Cells.Item(2, 9).Value2 = 0.09360395274144985  # I2: Ligand derived specificity of average expression value
$ws.Cells.Item(2, 10).Value2 = 0.09360395274144986  # J2: Ligand derived specificity of total expression value
$ws.Cells.Item(2, 13).Value2 = 133.7780026666667  # M2: Receptor average expression value
$ws.Cells.Item(2, 14).Value2 = 401.334008  # N2: Receptor total expression value
$ws.Cells.Item(2, 15).Value2 = 0.50863533211804  # O2: Receptor derived specificity of average expression value
$ws.Cells.Item(2, 16).Value2 = 0.5086353321180399  # P2: Receptor derived specificity of total expression value
$ws.Cells.Item(2, 17).Value2 = 246.7456776091769  # Q2: Edge average expression weight
$ws.Cells.Item(2, 18).Value2 = 2220.711098482592  # R2: Edge total expression weight
$ws.Cells.Item(2, 19).Value2 = 0.04761027759020867  # S2: Edge average expression derived specificity
$ws.Cells.Item(2, 20).Value2 = 0.04761027759020867  # T2: Edge total expression derived specificity

# Row 3
$ws.Cells.Item(3, 7).Value2 = 1.844441333333333  # G3: Ligand average expression value
$ws.Cells.Item(3, 8).Value2 = 5.533324  # H3: Ligand total expression value
$ws.Cells.Item(3, 9).Value2 = 0.09360395274144985  # I3: Ligand derived specificity of average expression value
$ws.Cells.Item(3, 10).Value2 = 0.09360395274144986  # J3: Ligand derived specificity of total expression value
$ws.Cells.Item(3, 15).Value2 = 0.1993888292903622  # O3: Receptor derived specificity of average expression value
$ws.Cells.Item(3, 16).Value2 = 0.1993888292903622  # P3: Receptor derived specificity of total expression value
$ws.Cells.Item(3, 17).Value2 = 96.72613891386801  # Q3: Edge average expression weight
$ws.Cells.Item(3, 18).Value2 = 870.5352502248121  # R3: Edge total expression weight
$ws.Cells.Item(3, 19).Value2 = 0.01866358255406808  # S3: Edge average expression derived specificity
$ws.Cells.Item(3, 20).Value2 = 0.01866358255406808  # T3: Edge total expression derived specificity

# Row 4
$ws.Cells.Item(4, 7).Value2 = 1.844441333333333  # G4: Ligand average expression value
$ws.Cells.Item(4, 8).Value2 = 5.533324  # H4: Ligand total expression value
$ws.Cells.Item(4, 9).Value2 = 0.09360395274144985  # I4: Ligand derived specificity of average expression value
$ws.Cells.Item(4, 10).Value2 = 0.09360395274144986  # J4: Ligand derived specificity of total expression value
$ws.Cells.Item(4, 13).Value2 = 21.197691  # M4: Receptor average expression value
$ws.Cells.Item(4, 14).Value2 = 63.593073  # N4: Receptor total expression value
$ws.Cells.Item(4, 15).Value2 = 0.08059542216956049  # O4: Receptor derived specificity of average expression value
$ws.Cells.Item(4, 16).Value2 = 0.08059542216956046  # P4: Receptor derived specificity of total expression value
$ws.Cells.Item(4, 17).Value2 = 39.09789745162801  # Q4: Edge average expression weight
$ws.Cells.Item(4, 18).Value2 = 351.8810770646521  # R4: Edge total expression weight
$ws.Cells.Item(4, 19).Value2 = 0.00754405008793674  # S4: Edge average expression derived specificity
$ws.Cells.Item(4, 20).Value2 = 0.007544050087936738  # T4: Edge total expression derived specificity

# Row 5
$ws.Cells.Item(5, 7).Value2 = 1.844441333333333  # G5: Ligand average expression value
$ws.Cells.Item(5, 8).Value2 = 5.533324  # H5: Ligand total expression value
$ws.Cells.Item(5, 9).Value2 = 0.09360395274144985  # I5: Ligand derived specificity of average expression value
$ws.Cells.Item(5, 10).Value2 = 0.09360395274144986  # J5: Ligand derived specificity of total expression value
$ws.Cells.Item(5, 13).Value2 = 55.59592133333333  # M5: Receptor average expression value
$ws.Cells.Item(5, 14).Value2 = 166.787764  # N5: Receptor total expression value
$ws.Cells.Item(5, 15).Value2 = 0.2113804164220374  # O5: Receptor derived specificity of average expression value
$ws.Cells.Item(5, 16).Value2 = 0.2113804164220373  # P5: Receptor derived specificity of total expression value
$ws.Cells.Item(5, 17).Value2 = 102.5434152719484  # Q5: Edge average expression weight
$ws.Cells.Item(5, 18).Value2 = 922.8907374475359  # R5: Edge total expression weight
$ws.Cells.Item(5, 19).Value2 = 0.01978604250923638  # S5: Edge average expression derived specificity
$ws.Cells.Item(5, 20).Value2 = 0.01978604250923637  # T5: Edge total expression derived specificity

# Row 6
$ws.Cells.Item(6, 9).Value2 = 0.3174745301946251  # I6: Ligand derived specificity of average expression value
$ws.Cells.Item(6, 10).Value2 = 0.3174745301946252  # J6: Ligand derived specificity of total expression value
$ws.Cells.Item(6, 13).Value2 = 133.7780026666667  # M6: Receptor average expression value
$ws.Cells.Item(6, 14).Value2 = 401.334008  # N6: Receptor total expression value
$ws.Cells.Item(6, 15).Value2 = 0.50863533211804  # O6: Receptor derived specificity of average expression value
$ws.Cells.Item(6, 16).Value2 = 0.5086353321180399  # P6: Receptor derived specificity of total expression value
$ws.Cells.Item(6, 17).Value2 = 836.882052330673  # Q6: Edge average expression weight
$ws.Cells.Item(6, 18).Value2 = 7531.938470976057  # R6: Edge total expression weight
$ws.Cells.Item(6, 19).Value2 = 0.1614787631045619  # S6: Edge average expression derived specificity
$ws.Cells.Item(6, 20).Value2 = 0.1614787631045619  # T6: Edge total expression derived specificity

# Row 7
$ws.Cells.Item(7, 9).Value2 = 0.3174745301946251  # I7: Ligand derived specificity of average expression value
$ws.Cells.Item(7, 10).Value2 = 0.3174745301946252  # J7: Ligand derived specificity of total expression value
$ws.Cells.Item(7, 15).Value2 = 0.1993888292903622  # O7: Receptor derived specificity of average expression value
$ws.Cells.Item(7, 16).Value2 = 0.1993888292903622  # P7: Receptor derived specificity of total expression value
$ws.Cells.Item(7, 18).Value2 = 2952.575842030642  # R7: Edge total expression weight
$ws.Cells.Item(7, 19).Value2 = 0.06330087490501406  # S7: Edge average expression derived specificity
$ws.Cells.Item(7, 20).Value2 = 0.06330087490501406  # T7: Edge total expression derived specificity

# Row 8
$ws.Cells.Item(8, 9).Value2 = 0.3174745301946251  # I8: Ligand derived specificity of average expression value
$ws.Cells.Item(8, 10).Value2 = 0.3174745301946252  # J8: Ligand derived specificity of total expression value
$ws.Cells.Item(8, 13).Value2 = 21.197691  # M8: Receptor average expression value
$ws.Cells.Item(8, 14).Value2 = 63.593073  # N8: Receptor total expression value
$ws.Cells.Item(8, 15).Value2 = 0.08059542216956049  # O8: Receptor derived specificity of average expression value
$ws.Cells.Item(8, 16).Value2 = 0.08059542216956046  # P8: Receptor derived specificity of total expression value
$ws.Cells.Item(8, 17).Value2 = 132.607504934529  # Q8: Edge average expression weight
$ws.Cells.Item(8, 18).Value2 = 1193.467544410761  # R8: Edge total expression weight
$ws.Cells.Item(8, 19).Value2 = 0.02558699378911869  # S8: Edge average expression derived specificity
$ws.Cells.Item(8, 20).Value2 = 0.02558699378911869  # T8: Edge total expression derived specificity

# Row 9
$ws.Cells.Item(9, 9).Value2 = 0.3174745301946251  # I9: Ligand derived specificity of average expression value
$ws.Cells.Item(9, 10).Value2 = 0.3174745301946252  # J9: Ligand derived specificity of total expression value
$ws.Cells.Item(9, 13).Value2 = 55.59592133333333  # M9: Receptor average expression value
$ws.Cells.Item(9, 14).Value2 = 166.787764  # N9: Receptor total expression value
$ws.Cells.Item(9, 15).Value2 = 0.2113804164220374  # O9: Receptor derived specificity of average expression value
$ws.Cells.Item(9, 16).Value2 = 0.2113804164220373  # P9: Receptor derived specificity of total expression value
$ws.Cells.Item(9, 17).Value2 = 347.7943146048164  # Q9: Edge average expression weight
$ws.Cells.Item(9, 18).Value2 = 3130.148831443348  # R9: Edge total expression weight
$ws.Cells.Item(9, 19).Value2 = 0.06710789839593054  # S9: Edge average expression derived specificity
$ws.Cells.Item(9, 20).Value2 = 0.06710789839593052  # T9: Edge total expression derived specificity

# Row 10
$ws.Cells.Item(10, 7).Value2 = 5.721023666666667  # G10: Ligand average expression value
$ws.Cells.Item(10, 8).Value2 = 17.163071  # H10: Ligand total expression value
$ws.Cells.Item(10, 9).Value2 = 0.2903374692647943  # I10: Ligand derived specificity of average expression value
$ws.Cells.Item(10, 10).Value2 = 0.2903374692647943  # J10: Ligand derived specificity of total expression value
$ws.Cells.Item(10, 13).Value2 = 133.7780026666667  # M10: Receptor average expression value
$ws.Cells.Item(10, 14).Value2 = 401.334008  # N10: Receptor total expression value
$ws.Cells.Item(10, 15).Value2 = 0.50863533211804  # O10: Receptor derived specificity of average expression value
$ws.Cells.Item(10, 16).Value2 = 0.5086353321180399  # P10: Receptor derived specificity of total expression value
$ws.Cells.Item(10, 17).Value2 = 765.3471193353965  # Q10: Edge average expression weight
$ws.Cells.Item(10, 18).Value2 = 6888.124074018569  # R10: Edge total expression weight
$ws.Cells.Item(10, 19).Value2 = 0.1476758951058099  # S10: Edge average expression derived specificity
$ws.Cells.Item(10, 20).Value2 = 0.1476758951058098  # T10: Edge total expression derived specificity

# Row 11
$ws.Cells.Item(11, 7).Value2 = 5.721023666666667  # G11: Ligand average expression value
$ws.Cells.Item(11, 8).Value2 = 17.163071  # H11: Ligand total expression value
$ws.Cells.Item(11, 9).Value2 = 0.2903374692647943  # I11: Ligand derived specificity of average expression value
$ws.Cells.Item(11, 10).Value2 = 0.2903374692647943  # J11: Ligand derived specificity of total expression value
$ws.Cells.Item(11, 15).Value2 = 0.1993888292903622  # O11: Receptor derived specificity of average expression value
$ws.Cells.Item(11, 16).Value2 = 0.1993888292903622  # P11: Receptor derived specificity of total expression value
$ws.Cells.Item(11, 17).Value2 = 300.021757217647  # Q11: Edge average expression weight
$ws.Cells.Item(11, 18).Value2 = 2700.195814958823  # R11: Edge total expression weight
$ws.Cells.Item(11, 19).Value2 = 0.05789004809583385  # S11: Edge average expression derived specificity
$ws.Cells.Item(11, 20).Value2 = 0.05789004809583384  # T11: Edge total expression derived specificity

# Row 12
$ws.Cells.Item(12, 7).Value2 = 5.721023666666667  # G12: Ligand average expression value
$ws.Cells.Item(12, 8).Value2 = 17.163071  # H12: Ligand total expression value
$ws.Cells.Item(12, 9).Value2 = 0.2903374692647943  # I12: Ligand derived specificity of average expression value
$ws.Cells.Item(12, 10).Value2 = 0.2903374692647943  # J12: Ligand derived specificity of total expression value
$ws.Cells.Item(12, 13).Value2 = 21.197691  # M12: Receptor average expression value
$ws.Cells.Item(12, 14).Value2 = 63.593073  # N12: Receptor total expression value
$ws.Cells.Item(12, 15).Value2 = 0.08059542216956049  # O12: Receptor derived specificity of average expression value
$ws.Cells.Item(12, 16).Value2 = 0.08059542216956046  # P12: Receptor derived specificity of total expression value
$ws.Cells.Item(12, 17).Value2 = 121.272491889687  # Q12: Edge average expression weight
$ws.Cells.Item(12, 18).Value2 = 1091.452427007183  # R12: Edge total expression weight
$ws.Cells.Item(12, 19).Value2 = 0.02339987090703789  # S12: Edge average expression derived specificity
$ws.Cells.Item(12, 20).Value2 = 0.02339987090703788  # T12: Edge total expression derived specificity

# Row 13
$ws.Cells.Item(13, 7).Value2 = 5.721023666666667  # G13: Ligand average expression value
$ws.Cells.Item(13, 8).Value2 = 17.163071  # H13: Ligand total expression value
$ws.Cells.Item(13, 9).Value2 = 0.2903374692647943  # I13: Ligand derived specificity of average expression value
$ws.Cells.Item(13, 10).Value2 = 0.2903374692647943  # J13: Ligand derived specificity of total expression value
$ws.Cells.Item(13, 13).Value2 = 55.59592133333333  # M13: Receptor average expression value
$ws.Cells.Item(13, 14).Value2 = 166.787764  # N13: Receptor total expression value
$ws.Cells.Item(13, 15).Value2 = 0.2113804164220374  # O13: Receptor derived specificity of average expression value
$ws.Cells.Item(13, 16).Value2 = 0.2113804164220373  # P13: Receptor derived specificity of total expression value
$ws.Cells.Item(13, 17).Value2 = 318.0655817181382  # Q13: Edge average expression weight
$ws.Cells.Item(13, 18).Value2 = 2862.590235463244  # R13: Edge total expression weight
$ws.Cells.Item(13, 19).Value2 = 0.06137165515611268  # S13: Edge average expression derived specificity
$ws.Cells.Item(13, 20).Value2 = 0.06137165515611267  # T13: Edge total expression derived specificity

# Row 14
$ws.Cells.Item(14, 7).Value2 = 5.883520333333333  # G14: Ligand average expression value
$ws.Cells.Item(14, 8).Value2 = 17.650561  # H14: Ligand total expression value
$ws.Cells.Item(14, 9).Value2 = 0.2985840477991308  # I14: Ligand derived specificity of average expression value
$ws.Cells.Item(14, 10).Value2 = 0.2985840477991308  # J14: Ligand derived specificity of total expression value
$ws.Cells.Item(14, 13).Value2 = 133.7780026666667  # M14: Receptor average expression value
$ws.Cells.Item(14, 14).Value2 = 401.334008  # N14: Receptor total expression value
$ws.Cells.Item(14, 15).Value2 = 0.50863533211804  # O14: Receptor derived specificity of average expression value
$ws.Cells.Item(14, 16).Value2 = 0.5086353321180399  # P14: Receptor derived specificity of total expression value
$ws.Cells.Item(14, 17).Value2 = 787.0855988420543  # Q14: Edge average expression weight
$ws.Cells.Item(14, 18).Value2 = 7083.770389578489  # R14: Edge total expression weight
$ws.Cells.Item(14, 19).Value2 = 0.1518703963174597  # S14: Edge average expression derived specificity
$ws.Cells.Item(14, 20).Value2 = 0.1518703963174596  # T14: Edge total expression derived specificity

# Row 15
$ws.Cells.Item(15, 7).Value2 = 5.883520333333333  # G15: Ligand average expression value
$ws.Cells.Item(15, 8).Value2 = 17.650561  # H15: Ligand total expression value
$ws.Cells.Item(15, 9).Value2 = 0.2985840477991308  # I15: Ligand derived specificity of average expression value
$ws.Cells.Item(15, 10).Value2 = 0.2985840477991308  # J15: Ligand derived specificity of total expression value
$ws.Cells.Item(15, 15).Value2 = 0.1993888292903622  # O15: Receptor derived specificity of average expression value
$ws.Cells.Item(15, 16).Value2 = 0.1993888292903622  # P15: Receptor derived specificity of total expression value
$ws.Cells.Item(15, 17).Value2 = 308.543402698577  # Q15: Edge average expression weight
$ws.Cells.Item(15, 18).Value2 = 2776.890624287193  # R15: Edge total expression weight
$ws.Cells.Item(15, 19).Value2 = 0.05953432373544625  # S15: Edge average expression derived specificity
$ws.Cells.Item(15, 20).Value2 = 0.05953432373544625  # T15: Edge total expression derived specificity

# Row 16
$ws.Cells.Item(16, 7).Value2 = 5.883520333333333  # G16: Ligand average expression value
$ws.Cells.Item(16, 8).Value2 = 17.650561  # H16: Ligand total expression value
$ws.Cells.Item(16, 9).Value2 = 0.2985840477991308  # I16: Ligand derived specificity of average expression value
$ws.Cells.Item(16, 10).Value2 = 0.2985840477991308  # J16: Ligand derived specificity of total expression value
$ws.Cells.Item(16, 13).Value2 = 21.197691  # M16: Receptor average expression value
$ws.Cells.Item(16, 14).Value2 = 63.593073  # N16: Receptor total expression value
$ws.Cells.Item(16, 15).Value2 = 0.08059542216956049  # O16: Receptor derived specificity of average expression value
$ws.Cells.Item(16, 16).Value2 = 0.08059542216956046  # P16: Receptor derived specificity of total expression value
$ws.Cells.Item(16, 17).Value2 = 124.717046018217  # Q16: Edge average expression weight
$ws.Cells.Item(16, 18).Value2 = 1122.453414163953  # R16: Edge total expression weight
$ws.Cells.Item(16, 19).Value2 = 0.02406450738546718  # S16: Edge average expression derived specificity
$ws.Cells.Item(16, 20).Value2 = 0.02406450738546717  # T16: Edge total expression derived specificity

# Row 17
$ws.Cells.Item(17, 7).Value2 = 5.883520333333333  # G17: Ligand average expression value
$ws.Cells.Item(17, 8).Value2 = 17.650561  # H17: Ligand total expression value
$ws.Cells.Item(17, 9).Value2 = 0.2985840477991308  # I17: Ligand derived specificity of average expression value
$ws.Cells.Item(17, 10).Value2 = 0.2985840477991308  # J17: Ligand derived specificity of total expression value
$ws.Cells.Item(17, 13).Value2 = 55.59592133333333  # M17: Receptor average expression value
$ws.Cells.Item(17, 14).Value2 = 166.787764  # N17: Receptor total expression value
$ws.Cells.Item(17, 15).Value2 = 0.2113804164220374  # O17: Receptor derived specificity of average expression value
$ws.Cells.Item(17, 16).Value2 = 0.2113804164220373  # P17: Receptor derived specificity of total expression value
$ws.Cells.Item(17, 17).Value2 = 327.0997336150671  # Q17: Edge average expression weight
$ws.Cells.Item(17, 18).Value2 = 2943.897602535604  # R17: Edge total expression weight
$ws.Cells.Item(17, 19).Value2 = 0.06311482036075779  # S17: Edge average expression derived specificity
$ws.Cells.Item(17, 20).Value2 = 0.06311482036075777  # T17: Edge total expression derived specificity
